$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") rows 2-289 all hold the serial date value 45203
# (2023-10-04) and must be bumped to 45204 (2023-10-05).
$ws.Range("C2:C289").Value = 45204
